# Generate Report for Archive
# Reorders rows 5-7 on the Overview / zh-cn / de-de sheets so that the
# ce8f7aa5 entry (which is now "In Translation") moves up to row 5, and the
# 51c1ed14 / 8567be4c entries shift down to rows 6 and 7 respectively.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet 1: Overview
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A5").Value = "ce8f7aa5-7343-4678-970f-1f723ba93e36.md"
$ws1.Range("B5").Value = "In Translation"
$ws1.Range("C5").Value = "In Translation"
$ws1.Range("D5").Value = "2016-32-21 06:32:42"

$ws1.Range("A6").Value = "51c1ed14-94c1-4234-a34c-84ccff7282fe.md"
$ws1.Range("B6").Value = "Ready for handoff"
$ws1.Range("C6").Value = "Ready for handoff"
$ws1.Range("D6").Value = "2016-30-21 06:30:36"

$ws1.Range("A7").Value = "8567be4c-b30c-46c3-85f4-62ca48e51d69.md"
$ws1.Range("B7").Value = "Ready for handoff"
$ws1.Range("C7").Value = "Ready for handoff"
$ws1.Range("D7").Value = "2016-33-21 06:33:06"

# ----------------------------------------------------------------------
# Sheet 2: zh-cn
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A5").Value = "ce8f7aa5-7343-4678-970f-1f723ba93e36.md"
$ws2.Range("C5").Value = "In Translation"
$ws2.Range("D5").Value = "ce8f7aa5-7343-4678-970f-1f723ba93e36.a4a358ac64abc54f25868d9795d3946ad6dbdbbb.zh-cn.xlf"
$ws2.Range("E5").Value = "2016-03-21 06:32:38"

$ws2.Range("A6").Value = "51c1ed14-94c1-4234-a34c-84ccff7282fe.md"
$ws2.Range("C6").Value = "Ready for handoff"
$ws2.Range("D6").Value = "51c1ed14-94c1-4234-a34c-84ccff7282fe.cd45b03d24c9259dce136154e9fae89f337ebee1.zh-cn.xlf"
$ws2.Range("E6").Value = "2016-03-21 06:30:33"

$ws2.Range("A7").Value = "8567be4c-b30c-46c3-85f4-62ca48e51d69.md"
$ws2.Range("C7").Value = "Ready for handoff"
$ws2.Range("D7").Value = "8567be4c-b30c-46c3-85f4-62ca48e51d69.e27cbadae0e4305f524b21969f3a05bb2e472570.zh-cn.xlf"
$ws2.Range("E7").Value = "2016-03-21 06:33:01"

# ----------------------------------------------------------------------
# Sheet 3: de-de
# ----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A5").Value = "ce8f7aa5-7343-4678-970f-1f723ba93e36.md"
$ws3.Range("C5").Value = "In Translation"
$ws3.Range("D5").Value = "ce8f7aa5-7343-4678-970f-1f723ba93e36.a4a358ac64abc54f25868d9795d3946ad6dbdbbb.de-de.xlf"
$ws3.Range("E5").Value = "2016-03-21 06:32:42"

$ws3.Range("A6").Value = "51c1ed14-94c1-4234-a34c-84ccff7282fe.md"
$ws3.Range("C6").Value = "Ready for handoff"
$ws3.Range("D6").Value = "51c1ed14-94c1-4234-a34c-84ccff7282fe.cd45b03d24c9259dce136154e9fae89f337ebee1.de-de.xlf"
$ws3.Range("E6").Value = "2016-03-21 06:30:36"

$ws3.Range("A7").Value = "8567be4c-b30c-46c3-85f4-62ca48e51d69.md"
$ws3.Range("C7").Value = "Ready for handoff"
$ws3.Range("D7").Value = "8567be4c-b30c-46c3-85f4-62ca48e51d69.e27cbadae0e4305f524b21969f3a05bb2e472570.de-de.xlf"
$ws3.Range("E7").Value = "2016-03-21 06:33:06"

Write-Host "Row reorder complete"
